$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @(-8.5027979459025147,-7.9005035341938283,-10.318909289399613,-17.065488032827023,-10.387075778692937,-7.8040668420551693,-8.611666673787381,-9.4126301369939078,-9.6600555311235734,-8.5039312075102167,-15.563337054149205,-8.4848063598336871,-12.137290610967435,-9.0547490243011541,-11.694034329364168,-9.8220101425661852,-8.5927406695942459,-9.9168666754148767,-10.984637576014704,-7.709814212591227,-7.337428974937434,-7.347942730991571,-8.7448267011052785,-8.1593341918530911,-9.2404798601130906,-7.7491814070169527,-8.4905371128445903,-12.819274273468421,-7.8533143801379994,-13.037404774714156,-10.714499405569468,-7.169921157353504,-11.975798081783259,-8.8914026735732499,-9.5252696271530439,-10.901006427368637,-9.1771084613994507,-10.064113850840329,-9.5254933030077265,-9.597385068813379,-9.56618542480971,-8.6568647354278134,-9.8862944619730531,-12.037671003491839,-8.077568918913471,-8.8502601036849278,-8.9733312742911995,-7.6208286827174199,-12.964057415181925,-8.8908861662947896,-12.160402915218267,-9.94154650321601,-8.6531966846780932,-10.173497377880368,-13.625364101004836,-11.429364422369915,-9.1834277332688963,-9.3662429501638247,-8.0208282415358507,-13.910305551863209,-7.136018365531549,-8.4272270768233515,-10.983037927536399,-9.0783476168401904,-12.248573880606244,-7.9424120992298812,-8.3590542072904039,-11.156768611735572,-8.1975140261739146,-10.671045867442468,-8.0482698522182989,-7.4325269519806296,-8.1915861032938633,-10.290314527176747,-9.3882999617092349,-9.8063375245711182,-8.6684769109030242,-8.5579118957599771,-13.760413663186975,-8.0743307413300709,-7.8618507396980846,-17.106835836772852,-7.7441934213394585,-6.8707472157951113,-8.3648655470182245,-8.7844540185149178,-13.748142954373503,-14.926683418740831,-9.1018937270407623,-10.664899842278389,-9.6413022081032675,-8.3066664838997326,-12.038628282745353,-8.8429855324554136,-8.0745461086976711,-9.6967706661314743,-8.4784011465475881,-8.1759850459632268,-12.352524315318298,-9.2533255093138891)
$row2 = @(-7.0544416896007984,-8.3723429635618309,-9.7578121545372962,-15.197018984563766,-8.8421868775263679,-7.2991155523763878,-8.1771748086211833,-8.9280921073815183,-9.1326893233370292,-8.0314468244242221,-14.792285522663523,-8.9628233366257888,-10.520447817632274,-8.510254150981174,-11.062080871808927,-9.3797912805248256,-8.1652376741880737,-9.3758787670459327,-10.508921548679,-6.3210829016517724,-7.8866935168505181,-7.86429920340266,-8.2462857957348827,-7.7048271843325704,-8.7381773865475676,-7.3388947966742482,-8.0526805559421106,-12.173367771768945,-8.2889673440528266,-12.398800368111779,-10.160965039623134,-7.6871147589545439,-10.341966154007562,-8.3938911601397628,-9.0324137738597958,-9.2798933906241512,-8.6729565639070323,-9.5339057312027542,-8.9920694058083761,-9.1098366066006502,-9.0368600430538883,-8.1116117047771699,-10.337057782278997,-11.355347243055995,-7.6393023728310521,-8.3253741979201692,-8.4571943040885582,-8.1886734258680907,-12.26275600477803,-9.3329845321471243,-10.570001664877536,-8.4194209921715437,-8.1489828583714594,-9.6333681823258051,-12.928172313343241,-10.797533579415967,-8.730956594258096,-8.8567925470474975,-7.5779634683864554,-13.270663474750124,-7.7559330593191227,-7.9517277047454904,-10.334467370063306,-7.5712210764763324,-11.564831785916507,-7.5319127144974987,-7.8430906147727848,-9.6281259732488245,-7.7437769162531502,-10.143728032285726,-7.5877485413337435,-8.0195016720958865,-6.695645773734908,-8.7639491167189263,-9.7756517839565422,-9.2357781503405629,-8.2898907102009094,-8.0790943045891819,-12.092473449857044,-7.6109915560507471,-7.4410328637207108,-16.236831111664749,-7.3160221061232606,-7.4337192166621699,-7.9062300112024184,-8.2521370556780909,-12.102817010451753,-14.239173559716194,-8.629129534454016,-9.1875184109107977,-9.1018369547112208,-7.80400145025295,-11.319324821199917,-8.3246287807546953,-7.5742574060898793,-9.1641568269430458,-8.0256967701764257,-7.7188132083926551,-11.732681854748863,-8.6795576212226653)

for ($c = 0; $c -lt 100; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $row1[$c]
    $ws.Cells.Item(2, $c + 1).Value = $row2[$c]
}
